$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 10

$ws.Cells.Item($row, 1).Value = 9.0
$ws.Cells.Item($row, 2).Value = "Monday, Jan 09"
$ws.Cells.Item($row, 3).Value = "2:30 PM"
$ws.Cells.Item($row, 4).Value = "LO3993"
$ws.Cells.Item($row, 5).Value = "Warsaw"
$ws.Cells.Item($row, 6).Value = "(WAW)"
$ws.Cells.Item($row, 7).Value = "LOT "
$ws.Cells.Item($row, 8).Value = "E170"
$ws.Cells.Item($row, 9).Value = "(SP-LDF)"
$ws.Cells.Item($row, 10).Value = "2:41 PM"
$ws.Cells.Item($row, 11).Borders.LineStyle = -4142
$ws.Cells.Item($row, 12).Value = "0 hours, 11 minutes"
$ws.Cells.Item($row, 13).Borders.LineStyle = -4142
